$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 26490
$ws.Range("J63").Value = 26490
$ws.Range("L63").Value = 26490
$ws.Range("N63").Value = -27738
$ws.Range("H66").Value = 26490
$ws.Range("J66").Value = 26490
$ws.Range("L66").Value = 79470
$ws.Range("N66").Value = -85710
$ws.Range("H125").Value = 788.93335
$ws.Range("I125").Value = 799.75
$ws.Range("J125").Value = 785
$ws.Range("K125").Value = 7197.75
$ws.Range("L125").Value = 7065
$ws.Range("M125").Value = -4737.75
$ws.Range("N125").Value = -11985
$ws.Range("H132").Value = 2337.1162
$ws.Range("I132").Value = 2345.1428
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7035.428400000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -4505.428400000001
$ws.Range("N132").Value = -11060
$ws.Range("H135").Value = 488.8846
$ws.Range("I135").Value = 459.5909
$ws.Range("J135").Value = 650
$ws.Range("K135").Value = 4136.3181
$ws.Range("L135").Value = 5850
$ws.Range("M135").Value = -1601.3181
$ws.Range("N135").Value = -10920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 924.6667
$ws.Range("I2").Value = 907.3333
$ws.Range("J2").Value = 942
$ws.Range("K2").Value = 907.3333
$ws.Range("L2").Value = 942
$ws.Range("M2").Value = -794.3333
$ws.Range("N2").Value = -1168
$ws.Range("H3").Value = 25007696
$ws.Range("I3").Value = 502.5
$ws.Range("J3").Value = 50014890
$ws.Range("K3").Value = 502.5
$ws.Range("L3").Value = 50014890
$ws.Range("M3").Value = -387.5
$ws.Range("N3").Value = -50015120
$ws.Range("H6").Value = 11666.167
$ws.Range("J6").Value = 11666.167
$ws.Range("L6").Value = 11666.167
$ws.Range("N6").Value = -12012.167
$ws.Range("H32").Value = 18186628
$ws.Range("I32").Value = 21741316
$ws.Range("K32").Value = 21741316
$ws.Range("M32").Value = -21741029
$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 50000
$ws.Range("N66").Value = -56864
$ws.Range("H116").Value = 924.6667
$ws.Range("I116").Value = 907.3333
$ws.Range("J116").Value = 942
$ws.Range("K116").Value = 907.3333
$ws.Range("L116").Value = 942
$ws.Range("M116").Value = 1386.6667
$ws.Range("N116").Value = -5530
$ws.Range("H132").Value = 2078.9092
$ws.Range("I132").Value = 1886
$ws.Range("J132").Value = 2492.2856
$ws.Range("K132").Value = 5658
$ws.Range("L132").Value = 7476.8568
$ws.Range("M132").Value = -3128
$ws.Range("N132").Value = -12536.8568
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 924.6667
$ws.Range("I3").Value = 907.3333
$ws.Range("J3").Value = 942
$ws.Range("K3").Value = 907.3333
$ws.Range("L3").Value = 942
$ws.Range("M3").Value = -793.3333
$ws.Range("N3").Value = -1170
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 9133.333000000001
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 13500
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 13500
$ws.Range("M12").Value = -230
$ws.Range("N12").Value = -13840
$ws.Range("H17").Value = 17699.6
$ws.Range("I17").Value = 498
$ws.Range("J17").Value = 22000
$ws.Range("K17").Value = 498
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = -324
$ws.Range("N17").Value = -22348
$ws.Range("H74").Value = 25016.666
$ws.Range("I74").Value = 27800
$ws.Range("J74").Value = 24460
$ws.Range("K74").Value = 27800
$ws.Range("L74").Value = 24460
$ws.Range("M74").Value = -26926
$ws.Range("N74").Value = -26208
$ws.Range("H77").Value = 25016.666
$ws.Range("I77").Value = 27800
$ws.Range("J77").Value = 24460
$ws.Range("K77").Value = 83400
$ws.Range("L77").Value = 73380
$ws.Range("M77").Value = -79032
$ws.Range("N77").Value = -82116
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 20000
$ws.Range("N88").Value = -20812
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("N91").Value = -22808
$ws.Range("H134").Value = 3248.6538
$ws.Range("I134").Value = 1267.125
$ws.Range("J134").Value = 6419.1
$ws.Range("K134").Value = 3801.375
$ws.Range("L134").Value = 19257.3
$ws.Range("M134").Value = -1266.375
$ws.Range("N134").Value = -24327.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1062.058
$ws.Range("I68").Value = 1013.5172
$ws.Range("J68").Value = 1097.25
$ws.Range("K68").Value = 3040.5516
$ws.Range("L68").Value = 3291.75
$ws.Range("M68").Value = -2229.5516
$ws.Range("N68").Value = -4913.75
$ws.Range("H71").Value = 1062.058
$ws.Range("I71").Value = 1013.5172
$ws.Range("J71").Value = 1097.25
$ws.Range("K71").Value = 9121.6548
$ws.Range("L71").Value = 9875.25
$ws.Range("M71").Value = -5065.6548
$ws.Range("N71").Value = -17987.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 8900
$ws.Range("J17").Value = 8900
$ws.Range("L17").Value = 8900
$ws.Range("N17").Value = -9236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1771.5385
$ws.Range("I82").Value = 1503
$ws.Range("J82").Value = 2666.6667
$ws.Range("K82").Value = 1503
$ws.Range("L82").Value = 2666.6667
$ws.Range("M82").Value = -1142
$ws.Range("N82").Value = -3388.6667
$ws.Range("H85").Value = 1771.5385
$ws.Range("I85").Value = 1503
$ws.Range("J85").Value = 2666.6667
$ws.Range("K85").Value = 1503
$ws.Range("L85").Value = 2666.6667
$ws.Range("M85").Value = -255
$ws.Range("N85").Value = -5162.6667
$ws.Range("H133").Value = 50292.668
$ws.Range("J133").Value = 50292.668
$ws.Range("L133").Value = 50292.668
$ws.Range("N133").Value = -55352.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1433400
$ws.Range("J7").Value = 6560
$ws.Range("L7").Value = 6560
$ws.Range("N7").Value = -6786
$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71248
$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -216240
